# Fruta / hortaliza, semanal
# Insert 4 new rows of weekly price data at the top of the "Papa" detail
# block (rows 525-528), pushing the existing rows 525-594 down to 529-598.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 525:598 down by inserting 4 blank rows at 525 (carries over
# formatting/styles from the row being pushed down, same as Excel's
# Insert Copied Cells / Insert Rows behaviour).
$ws.Range("A525:R528").Insert()

# --- New row 525 ---
$ws.Range("A525").Value = 9
$ws.Range("B525").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C525").Value = "Metropolitana"
$ws.Range("D525").Value2 = 44449
$ws.Range("E525").Value = 13
$ws.Range("F525").Value = 100114001
$ws.Range("G525").Value = "Papa"
$ws.Range("H525").Value = "Asterix"
$ws.Range("I525").Value = "1a (cosecha lavada)"
$ws.Range("J525").Value = 340
$ws.Range("K525").Value = 9000
$ws.Range("L525").Value = 10000
$ws.Range("M525").Value = 9500
$ws.Range("N525").Value = "`$/malla 25 kilos"
$ws.Range("O525").Value = "Provincia de Melipilla"
$ws.Range("P525").Value = 380
$ws.Range("Q525").Value = 25
$ws.Range("R525").Value = "Hortaliza"

# --- New row 526 ---
$ws.Range("A526").Value = 9
$ws.Range("B526").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C526").Value = "Metropolitana"
$ws.Range("D526").Value2 = 44449
$ws.Range("E526").Value = 13
$ws.Range("F526").Value = 100114001
$ws.Range("G526").Value = "Papa"
$ws.Range("H526").Value = "Asterix"
$ws.Range("I526").Value = "1a (cosecha)"
$ws.Range("J526").Value = 232
$ws.Range("K526").Value = 7000
$ws.Range("L526").Value = 8000
$ws.Range("M526").Value = 7500
$ws.Range("N526").Value = "`$/saco 25 kilos"
$ws.Range("O526").Value = "Provincia de Melipilla"
$ws.Range("P526").Value = 300
$ws.Range("Q526").Value = 25
$ws.Range("R526").Value = "Hortaliza"

# --- New row 527 ---
$ws.Range("A527").Value = 9
$ws.Range("B527").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C527").Value = "Metropolitana"
$ws.Range("D527").Value2 = 44449
$ws.Range("E527").Value = 13
$ws.Range("F527").Value = 100114001
$ws.Range("G527").Value = "Papa"
$ws.Range("H527").Value = "Rodeo"
$ws.Range("I527").Value = "1a (guarda lavada)"
$ws.Range("J527").Value = 430
$ws.Range("K527").Value = 8000
$ws.Range("L527").Value = 9000
$ws.Range("M527").Value = 8500
$ws.Range("N527").Value = "`$/malla 25 kilos"
$ws.Range("O527").Value = "Región de La Araucanía"
$ws.Range("P527").Value = 340
$ws.Range("Q527").Value = 25
$ws.Range("R527").Value = "Hortaliza"

# --- New row 528 ---
$ws.Range("A528").Value = 9
$ws.Range("B528").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C528").Value = "Metropolitana"
$ws.Range("D528").Value2 = 44449
$ws.Range("E528").Value = 13
$ws.Range("F528").Value = 100114001
$ws.Range("G528").Value = "Papa"
$ws.Range("H528").Value = "Rodeo"
$ws.Range("I528").Value = "1a (guarda)"
$ws.Range("J528").Value = 295
$ws.Range("K528").Value = 6000
$ws.Range("L528").Value = 7000
$ws.Range("M528").Value = 6498
$ws.Range("N528").Value = "`$/saco 25 kilos"
$ws.Range("O528").Value = "Región de La Araucanía"
$ws.Range("P528").Value = 260
$ws.Range("Q528").Value = 25
$ws.Range("R528").Value = "Hortaliza"
